$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "259.21"
    "D3"  = "26.95"
    "E3"  = "-1.01%"
    "D4"  = "4.668"
    "E4"  = "2.37%"
    "D5"  = "0.06027"
    "E5"  = "2.34%"
    "D6"  = "6.666"
    "E6"  = "0.49%"
    "D7"  = "0.8607"
    "E7"  = "0.38%"
    "D8"  = "0.9314"
    "E8"  = "0.58%"
    "D9"  = "0.1397"
    "E9"  = "-0.91%"
    "D10" = "0.04689"
    "E10" = "24.90%"
    "D11" = "0.07018"
    "E11" = "-0.82%"
    "D12" = "0.03115"
    "E12" = "-1.84%"
    "D13" = "0.09146"
    "E13" = "-0.36%"
    "D14" = "0.001552"
    "E14" = "0.03%"
    "D15" = "0.0006042"
    "E15" = "-0.44%"
    "D16" = "0.006099"
    "E16" = "0.16%"
    "D17" = "3.464"
    "E17" = "-1.49%"
    "D18" = "3.158"
    "E18" = "-1.23%"
    "E19" = "-1.71%"
    "D20" = "0.3113"
    "E20" = "0.22%"
    "D21" = "0.1296"
    "E21" = "1.48%"
    "D22" = "4.131"
    "E22" = "5.73%"
    "D23" = "0.04236"
    "E23" = "0.54%"
    "E24" = "-0.83%"
    "D25" = "0.004046"
    "E25" = "-5.91%"
    "E26" = "-0.17%"
    "E27" = "13.54%"
    "D40" = "0.03849"
    "E40" = "0.11%"
    "D41" = "0.1119"
    "E41" = "1.58%"
    "E42" = "-37.82%"
    "D43" = "0.002415"
    "E43" = "9.82%"
    "D44" = "0.01527"
    "E44" = "29.30%"
    "D45" = "0.00005105"
    "E45" = "-6.52%"
    "D46" = "0.00000000748"
    "E46" = "-0.18%"
    "E48" = "20.16%"
    "D49" = "0.00002096"
    "E49" = "-0.18%"
    "D50" = "0.0001996"
    "E50" = "-0.18%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force the cell to stay a plain text value (matching the original
    # inlineStr cells) instead of being auto-coerced into a number/percent.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
